# "added formula/calculation with tests and documentation"
#
# The XLSForm "survey" sheet gets a new "calculate" row that evaluates an
# expression combining the three existing boolean fields, plus a new
# "calculation" column (E) to hold that expression.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a blank row at position 7 (the old "end group" row and everything
# below it shifts from rows 7-14 down to rows 8-15).
$ws.Rows(7).Insert()

# New survey row: a "calculate" field named "calc".
$ws.Range("A7").Value = "calculate"
$ws.Range("B7").Value = "calc"

# New "calculation" header in column E (written before the row-7 cells
# below so new shared strings are appended in the same order as the
# original edit: calculate, calc, calculation, "All true:", the formula).
$ws.Range("E1").Value = "calculation"

$ws.Range("C7").Value = "All true:"
$ws.Range("E7").Value = '${show_slide} and ${show_group} and ${show_field}'

# Widen column E so the long calculation expression fits.
$ws.Columns("E").ColumnWidth = 48.5703125

# The sheet was left with cell G10 selected.
[void]$ws.Range("G10").Select()
